$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the existing header row, shifting all
# existing rows (header + data) down by one.
$ws.Rows("1:1").Insert()

# Write the new CSV-import header into the freshly inserted row 1.
$ws.Range("A1").Value = "col1"
$ws.Range("B1").Value = "col2"
$ws.Range("C1").Value = "col3"

# Row 2 (the old header "time"/"weight"/"radius") inherited the bold,
# centered, bordered header style when the row shifted down. Move that
# formatting onto the new row 1 instead, then strip it from row 2 so the
# old header reads as plain data, matching what the CSV-import GUI wrote.
$ws.Range("A2:C2").Copy($ws.Range("A1:C1"))
$ws.Range("A1").Value = "col1"
$ws.Range("B1").Value = "col2"
$ws.Range("C1").Value = "col3"
$ws.Range("A2:C2").ClearFormats()
